$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.447.64"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "1.733.06"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "322.76"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "0.4536"
$ws.Range("E7").Value = "  +6.72%  "
$ws.Range("D8").Value = "0.3531"
$ws.Range("E8").Value = "  -2.31%  "
$ws.Range("D9").Value = "0.07374"
$ws.Range("E9").Value = "  -2.59%  "
$ws.Range("D10").Value = "41.34"
$ws.Range("E10").Value = "  -2.89%  "
$ws.Range("D11").Value = "1.076"
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "20.41"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "5.926"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").Value = "7.059"
$ws.Range("E15").Value = "  -2.81%  "
$ws.Range("D16").Value = "1.726.96"
$ws.Range("E16").Value = "  -2.86%  "
$ws.Range("D17").Value = "91.01"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "'0.00001050"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").Value = "0.06327"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "16.62"
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("D22").Value = "5.725"
$ws.Range("E22").Value = "  -3.03%  "
$ws.Range("D23").Value = "27.453.33"
$ws.Range("E23").Value = "  -1.38%  "
$ws.Range("D24").Value = "11.07"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").Value = "2.059"
$ws.Range("E25").Value = "  -2.10%  "
$ws.Range("D26").Value = "161.87"
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("D27").Value = "19.87"
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("D28").Value = "1.925.31"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("D29").Value = "2.042"
$ws.Range("E29").Value = "  -4.36%  "
$ws.Range("D30").Value = "124.55"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").Value = "1.043"
$ws.Range("E31").Value = "  -6.46%  "
$ws.Range("D32").Value = "0.09138"
$ws.Range("E32").Value = "  +2.85%  "
$ws.Range("D33").Value = "3.649"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "5.356"
$ws.Range("E34").Value = "  -3.96%  "
$ws.Range("D35").Value = "0.02261"
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("D36").Value = "11.59"
$ws.Range("E36").Value = "  -5.24%  "
$ws.Range("D37").Value = "0.05951"
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("E38").Value = "  -2.82%  "
$ws.Range("D39").Value = "0.6216"
$ws.Range("E39").Value = "  -2.03%  "
$ws.Range("D40").Value = "4.857"
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("D41").Value = "1.187"
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("D42").Value = "1.368"
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("D43").Value = "7.714"
$ws.Range("E44").Value = "  -3.11%  "
$ws.Range("D45").Value = "3.696"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "0.5788"
$ws.Range("E46").Value = "  -1.46%  "
$ws.Range("D47").Value = "121.87"
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("D48").Value = "1.917"
$ws.Range("E48").Value = "  -3.35%  "
$ws.Range("D49").Value = "0.06831"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "1.109"
$ws.Range("E50").Value = "  -5.83%  "
$ws.Range("D51").Value = "70.86"
$ws.Range("E51").Value = "  -3.64%  "
